$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the spelling mistake in the "Remove Phasing Trace" related-shortcuts cell
$ws.Range("C48").Value = "PC, PV, TA"

# Remove the duplicate "Store Zoom" / "Restore Zoom" rows (rows 117:118),
# which duplicated the "Store Zoom setting" / "Restore stored Zoom setting" rows above.
# Everything below shifts up by two rows as a result.
$ws.Rows("117:118").Delete()

# Update the view to match the saved state
$ws.Application.ActiveWindow.ScrollRow = 36
$ws.Range("C49").Select()

$wb.Worksheets.Item("Sheet1").PageSetup.PrintArea = '$A$1:$C$126'
